# Finds the 1-based index of the paragraph whose text starts with $prefix.
function Find-ParaIndex($doc, $prefix) {
    $idx = 0
    $found = -1
    foreach ($p in $doc.Paragraphs) {
        $idx = $idx + 1
        if ($found -eq -1) {
            if ($p.Range.Text.StartsWith($prefix)) {
                $found = $idx
            }
        }
    }
    return $found
}

# Appends a new chunk of text at the very end of paragraph number $paraIndex
# (i.e. just before its paragraph mark), re-using the character formatting
# captured in $fmt (a FormattedText blob taken from an existing,
# already-correctly-formatted run). This lets the new text pick up
# sz / szCs / lang (and any other run formatting) exactly like the
# surrounding text, instead of being inserted with no formatting at all
# (which is what plain InsertAfter / Font.* property assignment does in
# this runtime).
function Append-Run($doc, $paraIndex, $fmt, $text) {
    $para = $doc.Paragraphs($paraIndex)
    $pr = $para.Range
    $insertPoint = $doc.Range($pr.Start, $pr.End - 1)
    $insertPoint.Collapse(0)
    $startPos = $insertPoint.End
    $insertPoint.FormattedText = $fmt
    $para2 = $doc.Paragraphs($paraIndex)
    $pr2 = $para2.Range
    $endPos = $doc.Range($pr2.Start, $pr2.End - 1).End
    $newRange = $doc.Range($startPos, $endPos)
    $newRange.Text = $text
}

$d = $word.ActiveDocument

# --- "Création d'une grille de jeu" + tab --------------------------------
# Grab formatting from the existing "Création d'une grille de jeu" run so
# the appended text shares sz=32 / szCs=32 / lang=fr-CH.
$i1 = Find-ParaIndex $d "Création d’une grille de jeu"
$r1 = $d.Paragraphs($i1).Range
$fmt1 = $d.Range($r1.Start, $r1.End - 2).FormattedText

Append-Run $d $i1 $fmt1 "en fonction de la valeur attribuée par le joueur au début"
Append-Run $d $i1 $fmt1 " et un tableau qui correspond à la taille de la grille"

# --- "Positionnement des bateaux ..." ------------------------------------
$i2 = Find-ParaIndex $d "Positionnement des bateaux"
$r2 = $d.Paragraphs($i2).Range
$fmt2 = $d.Range($r2.Start, $r2.End - 1).FormattedText

Append-Run $d $i2 $fmt2 ". En choisissant la case "
Append-Run $d $i2 $fmt2 "avec les touches "
Append-Run $d $i2 $fmt2 "et ensuite choisir l’orientation du bateau"
Append-Run $d $i2 $fmt2 " avec les flèches"

# --- "Le fait d’attaquer des cases du tableau ..." -----------------------
$i3 = Find-ParaIndex $d "Le fait d’attaquer"
$r3 = $d.Paragraphs($i3).Range
$fmt3 = $d.Range($r3.Start, $r3.End - 1).FormattedText

Append-Run $d $i3 $fmt3 " et si un joueur tire 30 coups dans l’eau il perd la partie"
Append-Run $d $i3 $fmt3 ". En choisissant la case "
Append-Run $d $i3 $fmt3 "sur laquelle nous tirons "
Append-Run $d $i3 $fmt3 "avec les touches"

# --- "Faire une condition de victoire ..." -------------------------------
# "touchés" -> "détruits entièrement"
$d.Content.Find.Execute("Faire une condition de victoire si tous les bateaux d’un joueur ont été touchés", $true, $false, $false, $false, $false, $true, 1, $false, "Faire une condition de victoire si tous les bateaux d’un joueur ont été détruits entièrement", 2)

"done"
